# #5: insurance, claim, debt, investment done
# Extend the "保險" (insurance) worksheet (sheet 6) from the partial
# A:D layout to the full A:K schema used by the other property sheets,
# and fix up the row-1 header labels (which were erroneously holding
# sample data instead of the column/field names).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# ---- Row 1: header / field-name labels ------------------------------
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"

# E1:K1 are brand-new header cells - pick up the same bold/centered/
# bordered look the existing B1:D1 header cells use before stamping
# in the field-name text.
$ws.Range("B1").Copy($ws.Range("E1:K1"))
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# ---- Data rows 2-8: fill in columns E:K with the common metadata ----
# property_category / category / date / legislator_name / legislator_id /
# source_file / index are the same shape every other sheet carries.
$indices = @(98, 99, 100, 101, 102, 103, 104)

for ($i = 0; $i -lt $indices.Length; $i++) {
    $row = $i + 2
    $idx = $indices[$i]

    $ws.Cells.Item($row, 5).Value = "insurance"     # E: property_category
    $ws.Cells.Item($row, 6).Value = "normal"         # F: category

    # G: date - format as Text first so the "2012-04-24" literal isn't
    # auto-converted into a date serial number by Excel's input parser.
    $dateCell = $ws.Cells.Item($row, 7)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2012-04-24"

    $ws.Cells.Item($row, 8).Value = "段宜康"          # H: legislator_name
    $ws.Cells.Item($row, 9).Value = 917              # I: legislator_id
    $ws.Cells.Item($row, 10).Value = "tmp25ce1"      # J: source_file
    $ws.Cells.Item($row, 11).Value = $idx            # K: index
}
